# Team member sheet: fill in the "Job Title" column (E) for the remaining
# rows with "member", and update the saved selection to F10.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rows 4-7 ("test2".."test5") had no Job Title set yet - give them "member"
# (rows 2 & 3 already have "CEO" / "CTO").
$ws.Range("E4").Value = "member"
$ws.Range("E5").Value = "member"
$ws.Range("E6").Value = "member"
$ws.Range("E7").Value = "member"

# Move/save the active selection as it was left in the edited workbook.
$ws.Range("F10").Select() | Out-Null
